# After splitting the test train method into multiple smaller modules.
# This populates the Predicted_next_Day_Price (AB), updates
# Predicted_Signal (AC) and Actual_Return (AD) columns for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AB column (Predicted_next_Day_Price) - newly added, all zero for rows 2-7
$ws.Range("AB2").Value = 0
$ws.Range("AB3").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AB7").Value = 0

# AC column (Predicted_Signal) updates
$ws.Range("AC3").Value = 1
$ws.Range("AC4").Value = 1
$ws.Range("AC7").Value = 1

# AD column (Actual_Return) updates
$ws.Range("AD4").Value = -0.002178269582643555
$ws.Range("AD5").Value = -0.006505413901501833
